$wb = $excel.ActiveWorkbook

# The "subgenus" field was removed from the Materials term list (per the
# third review round) - delete its entire column (header "subgenus" /
# value "${subgenus}") from the Materials sheet, shifting later columns left.
$ws = $wb.Worksheets.Item("Materials")
$ws.Range("AS1:AS2").EntireColumn.Delete()
